# Update for Vehicle Renewal
# Target sheet: vehicleRenewalTestData (4th sheet / sheetId 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vehicleRenewalTestData")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Prepare rows 3 and 4 by copying the formatting (styles, row height) of
#    row 2, so the new rows look the same as the existing data row.
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A2:I2").Copy()
$ws.Range("A4:I4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(2).RowHeight

# ---------------------------------------------------------------------------
# 2. Update the cell values. The order below intentionally introduces the
#    brand-new text values (FALSE, 4000, 25, TC 002, 400, TC 003, 401) in
#    the same sequence they are first needed so that the shared-strings
#    table is built up in the expected order.
# ---------------------------------------------------------------------------

# Row 2: TC 001 / 3000 / 4000 / VCL_ID_3 / false / PASSED / 25 / FALSE / False
$ws.Range("H2").Value = "'FALSE"
$ws.Range("C2").Value = "'4000"
$ws.Range("G2").Value = "'25"
$ws.Range("B2").Value = "'3000"
$ws.Range("I2").Value = "'False"

# Row 3: TC 002 / 3000 / 4000 / VCL_ID_3 / false / PASSED / 400 / FALSE / False
$ws.Range("A3").Value = "'TC 002"
$ws.Range("G3").Value = "'400"
$ws.Range("B3").Value = "'3000"
$ws.Range("C3").Value = "'4000"
$ws.Range("D3").Value = "'VCL_ID_3"
$ws.Range("E3").Value = "'false"
$ws.Range("F3").Value = "'PASSED"
$ws.Range("H3").Value = "'FALSE"
$ws.Range("I3").Value = "'False"

# Row 4: TC 003 / 1000 / 3000 / VCL_ID_3 / false / PASSED / 401 / FALSE / True
$ws.Range("A4").Value = "'TC 003"
$ws.Range("G4").Value = "'401"
$ws.Range("B4").Value = "'1000"
$ws.Range("C4").Value = "'3000"
$ws.Range("D4").Value = "'VCL_ID_3"
$ws.Range("E4").Value = "'false"
$ws.Range("F4").Value = "'PASSED"
$ws.Range("H4").Value = "'FALSE"
$ws.Range("I4").Value = "'True"

# ---------------------------------------------------------------------------
# 3. Extend the data validation lists from the single row (row 2) down to
#    cover rows 2:4, preserving the same validation settings/order.
# ---------------------------------------------------------------------------
$ws.Range("I2").Validation.Delete()
$ws.Range("I2:I4").Validation.Add(3, 1, 1, """True,False""")

$ws.Range("H2").Validation.Delete()
$ws.Range("H2:H4").Validation.Add(3, 1, 1, """TRUE,FALSE""")

$ws.Range("F2").Validation.Delete()
$ws.Range("F2:F4").Validation.Add(3, 1, 1, """PASSED,FAILED""")

$ws.Range("E2").Validation.Delete()
$ws.Range("E2:E4").Validation.Add(3, 1, 1, """true,false""")

# ---------------------------------------------------------------------------
# 4. Update the selected cell shown in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("H8").Select()
